$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; G=38.04655766666666; H=114.139673; I=0.8090698722086991; J=0.8090698722086992; K=3; M=1.707786666666667; N=5.12336; O=0.8764025646701329; P=0.8764025646701328; Q=64.97540389569777; R=584.7786350612801; S=0.7090709110010406; T=0.7090709110010406 }
    3 = @{ E=3; G=38.04655766666666; H=114.139673; I=0.8090698722086991; J=0.8090698722086992; K=3; M=0.240846; N=0.7225379999999999; O=0.1235974353298672; P=0.1235974353298672; Q=9.163361227785998; R=82.47025105007398; S=0.0999989612076586; T=0.0999989612076586 }
    4 = @{ E=3; G=1.617245333333334; H=4.851736000000001; I=0.03439113957782537; J=0.03439113957782537; K=3; M=1.707786666666667; N=5.12336; O=0.8764025646701329; P=0.8764025646701328; Q=2.761910016995556; R=24.85719015296; S=0.03014048292793467; T=0.03014048292793466 }
    5 = @{ E=3; G=1.617245333333334; H=4.851736000000001; I=0.03439113957782537; J=0.03439113957782537; K=3; M=0.240846; N=0.7225379999999999; O=0.1235974353298672; P=0.1235974353298672; Q=0.389507069552; R=3.505563625968; S=0.004250656649890707; T=0.004250656649890707 }
    6 = @{ E=3; G=7.361255; H=22.083765; I=0.1565389882134754; J=0.1565389882134754; K=3; M=1.707786666666667; N=5.12336; O=0.8764025646701329; P=0.8764025646701328; Q=12.57145313893333; R=113.1430782504; S=0.1371911707411576; T=0.1371911707411576 }
    7 = @{ E=3; G=7.361255; H=22.083765; I=0.1565389882134754; J=0.1565389882134754; K=3; M=0.240846; N=0.7225379999999999; O=0.1235974353298672; P=0.1235974353298672; Q=1.77292882173; R=15.95635939557; S=0.01934781747231788; T=0.01934781747231787 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
